$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A147").Value = "03_21/22"
$ws.Range("B147").Value = 155.2
$ws.Range("C147").Value = 7.1

$ws.Range("A148").Value = "04_21/22"
$ws.Range("B148").Value = 158.4
$ws.Range("C148").Value = 6.9
